# temp_exvel doesnt launch if no asc file selected
# Reorder the point rows: STN751, MENA and CTL03 move up earlier in the list
# (STN751 before STNRO, MENA before CTL4, CTL03 before STN13), and the
# trailing MENA row loses its ELEVATION value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2: STN751
Set-TextValue $ws.Range("A2") "STN751"
Set-TextValue $ws.Range("B2") "286440.7476"
Set-TextValue $ws.Range("C2") "6215456.2261"
Set-TextValue $ws.Range("D2") "114.0311"

# Row 3: STNRO
Set-TextValue $ws.Range("A3") "STNRO"
Set-TextValue $ws.Range("B3") "285930.9534"
Set-TextValue $ws.Range("C3") "6215397.8550"
Set-TextValue $ws.Range("D3") "115.5244"

# Row 4: STN03
Set-TextValue $ws.Range("A4") "STN03"
Set-TextValue $ws.Range("B4") "285968.9539"
Set-TextValue $ws.Range("C4") "6215310.3795"
Set-TextValue $ws.Range("D4") "117.3814"

# Row 5: CTL05
Set-TextValue $ws.Range("A5") "CTL05"
Set-TextValue $ws.Range("B5") "287047.0456"
Set-TextValue $ws.Range("C5") "6215313.2060"
Set-TextValue $ws.Range("D5") "133.2873"

# Row 6: MENA (no elevation)
Set-TextValue $ws.Range("A6") "MENA"
Set-TextValue $ws.Range("B6") "291928.7360"
Set-TextValue $ws.Range("C6") "6221563.0740"
$ws.Range("D6").ClearContents()

# Row 7: CTL4
Set-TextValue $ws.Range("A7") "CTL4"
Set-TextValue $ws.Range("B7") "286848.0783"
Set-TextValue $ws.Range("C7") "6215375.9164"
Set-TextValue $ws.Range("D7") "125.1082"

# Row 8: CTL03
Set-TextValue $ws.Range("A8") "CTL03"
Set-TextValue $ws.Range("B8") "286254.9359"
Set-TextValue $ws.Range("C8") "6215385.8195"
Set-TextValue $ws.Range("D8") "117.4396"

# Row 9: STN13
Set-TextValue $ws.Range("A9") "STN13"
Set-TextValue $ws.Range("B9") "287771.7689"
Set-TextValue $ws.Range("C9") "6215279.6768"
Set-TextValue $ws.Range("D9") "134.7438"

# Rows 10 (STN741) and 11 (STN12) are unchanged.
